{"js": "// Replace each \"a\u00d7b=\" multiplication expression with its updated version.\n// All source strings are unique within the document, so a direct\n// search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"11\u00d729=\", \"90\u00d750=\"],\n  [\"15\u00d715=\", \"60\u00d769=\"],\n  [\"77\u00d784=\", \"95\u00d765=\"],\n  [\"50\u00d763=\", \"25\u00d787=\"],\n  [\"14\u00d730=\", \"62\u00d751=\"],\n  [\"11\u00d797=\", \"96\u00d761=\"],\n  [\"71\u00d778=\", \"25\u00d715=\"],\n  [\"79\u00d753=\", \"90\u00d798=\"],\n  [\"37\u00d788=\", \"73\u00d724=\"],\n  [\"93\u00d766=\", \"75\u00d724=\"],\n  [\"32\u00d743=\", \"50\u00d773=\"],\n  [\"13\u00d738=\", \"27\u00d765=\"],\n  [\"65\u00d790=\", \"65\u00d728=\"],\n  [\"37\u00d768=\", \"11\u00d743=\"],\n  [\"70\u00d777=\", \"19\u00d768=\"],\n  [\"84\u00d743=\", \"67\u00d749=\"],\n  [\"98\u00d784=\", \"73\u00d743=\"],\n  [\"20\u00d717=\", \"86\u00d778=\"],\n  [\"94\u00d716=\", \"21\u00d771=\"],\n  [\"51\u00d771=\", \"31\u00d756=\"],\n  [\"59\u00d730=\", \"33\u00d773=\"],\n  [\"40\u00d739=\", \"45\u00d778=\"],\n  [\"29\u00d771=\", \"97\u00d791=\"],\n  [\"51\u00d737=\", \"47\u00d726=\"],\n  [\"82\u00d757=\", \"60\u00d736=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    continue;\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"a\u00d7b=\" multiplication expression with its updated version.\n# All source strings are unique within the document, so Find/Replace per\n# pair is unambiguous and safe to run across the whole document.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"11\u00d729=\", \"90\u00d750=\"),\n  @(\"15\u00d715=\", \"60\u00d769=\"),\n  @(\"77\u00d784=\", \"95\u00d765=\"),\n  @(\"50\u00d763=\", \"25\u00d787=\"),\n  @(\"14\u00d730=\", \"62\u00d751=\"),\n  @(\"11\u00d797=\", \"96\u00d761=\"),\n  @(\"71\u00d778=\", \"25\u00d715=\"),\n  @(\"79\u00d753=\", \"90\u00d798=\"),\n  @(\"37\u00d788=\", \"73\u00d724=\"),\n  @(\"93\u00d766=\", \"75\u00d724=\"),\n  @(\"32\u00d743=\", \"50\u00d773=\"),\n  @(\"13\u00d738=\", \"27\u00d765=\"),\n  @(\"65\u00d790=\", \"65\u00d728=\"),\n  @(\"37\u00d768=\", \"11\u00d743=\"),\n  @(\"70\u00d777=\", \"19\u00d768=\"),\n  @(\"84\u00d743=\", \"67\u00d749=\"),\n  @(\"98\u00d784=\", \"73\u00d743=\"),\n  @(\"20\u00d717=\", \"86\u00d778=\"),\n  @(\"94\u00d716=\", \"21\u00d771=\"),\n  @(\"51\u00d771=\", \"31\u00d756=\"),\n  @(\"59\u00d730=\", \"33\u00d773=\"),\n  @(\"40\u00d739=\", \"45\u00d778=\"),\n  @(\"29\u00d771=\", \"97\u00d791=\"),\n  @(\"51\u00d737=\", \"47\u00d726=\"),\n  @(\"82\u00d757=\", \"60\u00d736=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1  # wdFindContinue\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
